$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-8 with new values, keep E column labels (same except row 8)
$ws.Range("A2").Value = 52.0370991230011
$ws.Range("B2").Value = 50.3656177520752
$ws.Range("C2").Value = 8.682495832443237
$ws.Range("D2").Value = 17.76827049255371

$ws.Range("A3").Value = 135.5751173496246
$ws.Range("B3").Value = 199.3462336063385
$ws.Range("C3").Value = 10.06976366043091
$ws.Range("D3").Value = 18.02263808250427

$ws.Range("A4").Value = 38.15235161781311
$ws.Range("B4").Value = 61.09102964401245
$ws.Range("C4").Value = 7.977405786514282
$ws.Range("D4").Value = 15.25071573257446

$ws.Range("A5").Value = 84.50184178352356
$ws.Range("B5").Value = 32.64775466918945
$ws.Range("C5").Value = 5.094048023223877
$ws.Range("D5").Value = 16.10633516311646

$ws.Range("A6").Value = 84.64412069320679
$ws.Range("B6").Value = 66.65805077552795
$ws.Range("C6").Value = 16.42977380752563
$ws.Range("D6").Value = 672.990181684494

$ws.Range("A7").Value = 25.34094619750977
$ws.Range("B7").Value = 45.22903490066528
$ws.Range("C7").Value = 6.02359676361084
$ws.Range("D7").Value = 294.3282012939453

$ws.Range("A8").Value = 311.8839828968048
$ws.Range("B8").Value = 332.9332964420319
$ws.Range("C8").Value = 157.5405130386353
$ws.Range("D8").Value = 270.7927403450012
$ws.Range("E8").Value = "IMDB reviews"

# Delete row 9 entirely (previously held IMDB reviews row, now removed since data shifted up)
$ws.Rows.Item(9).Delete()
